# Generate Report for Archive
#
# 1. The Status column text "Ready for handoff" is now "In Translation"
#    (applies to every sheet: Overview, zh-cn, de-de).
# 2. Because the new status text is shorter, the Status column(s) on each
#    sheet are narrower than before (Overview!E:F, zh-cn!C, de-de!C).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Text -eq "Ready for handoff") {
                $cell.Value = "In Translation"
            }
        }
    }
}

# Narrow the Status columns to match the new (shorter) content.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").EntireColumn.ColumnWidth = 12.5
$overview.Range("F1").EntireColumn.ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
